$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67
$ws.Range("U2").Value = 1.92
$ws.Range("V2").Value = 1.77

# Row 3
$ws.Range("G3").Value = 2.38
$ws.Range("I3").Value = 3.6
$ws.Range("K3").Value = 1.8
$ws.Range("L3").Value = 4.5
$ws.Range("O3").Value = 1.8
$ws.Range("P3").Value = 1.91
$ws.Range("Q3").Value = 3.6
$ws.Range("R3").Value = 1.29
$ws.Range("S3").Value = 1.8
$ws.Range("T3").Value = 2
$ws.Range("W3").Value = 5
$ws.Range("AC3").Value = 4.75
$ws.Range("AM3").Value = 51
$ws.Range("AT3").Value = 2
$ws.Range("BA3").Value = 151

# Row 4
$ws.Range("BB4").Value = 500

# Row 5
$ws.Range("G5").Value = 4.33
$ws.Range("H5").Value = 2.9
$ws.Range("I5").Value = 2.05
$ws.Range("X5").Value = 19
